# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The employee's 24 "Periodo Mora" rows (16-39) are re-based from descending
# (2003 .. 1804) to ascending chronological order (1804 .. 2003), and the
# "Valor Mora" / "Salario Basico" figures are refreshed to the new database
# values: Salario Basico becomes a uniform 828116, and Valor Mora follows the
# period (oldest 9 periods -> 23437, next 6 -> 31249, newest 9 -> 33125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @("1804","1805","1806","1807","1808","1809","1810","1811","1812", `
             "1901","1902","1903","1904","1905","1906", `
             "1907","1908","1909","1910","1911","1912","2001","2002","2003")

$valorMora = @(23437,23437,23437,23437,23437,23437,23437,23437,23437, `
               31249,31249,31249,31249,31249,31249, `
               33125,33125,33125,33125,33125,33125,33125,33125,33125)

$salarioBasico = 828116

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico
}
